# ---------------------------------------------------------------------------
# multipumping.xlsx update:
#   - title (A3) gets a "(Stream Frequency=100MHz)" suffix
#   - H4/I4 header labels swap ("Freq per tick" <-> "Freq per elem")
#   - the G:J measurement columns for every data row (6-14) get refreshed
#     numbers from the finished "window summation" experiment run
#   - G gets a 6-decimal number format, H:J get a 2-decimal number format
#   - selection moves to K18
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Title: append " (Stream Frequency=100MHz)" after "... C value" -----
$cell = $ws.Range("A3")
$full = $cell.Characters().Text
$len = $full.Length
# Replace just the trailing " value" run with " value (Stream Frequency=100MHz)"
# so the preceding bold "W"/"C" runs are left alone.
$cell.Characters($len - 5, 6).Text = " value (Stream Frequency=100MHz)"

# Re-apply bold styling to the single-letter "W" and "C" runs, since setting
# .Text on the Characters range collapses run formatting.
$full2 = $cell.Characters().Text
$idxW = $full2.IndexOf("different ") + 10
$idxC = $full2.IndexOf(" and ") + 5
$ws.Range("A3").Characters($idxW + 1, 1).Font.Bold = $true
$ws.Range("A3").Characters($idxC + 1, 1).Font.Bold = $true

# --- 2. Swap the H4 / I4 header text ---------------------------------------
$ws.Range("H4").Value = "Freq per elem`n(MHz)"
$ws.Range("I4").Value = "Freq per tick`n(MHz)"

# --- 3. Refresh measurement data (columns G,H,I,J) for rows 6-14 -----------
# Each entry: row, G (time us), H (freq per elem), I (freq per tick), J (bandwidth)
$data = @(
    , @(6,  0.0027750000000000001, 90.1,  360.4,  1441.61)
    , @(7,  0.0027130000000000001, 46.08, 368.65, 1474.61)
    , @(8,  0.002686,              23.27, 372.29, 1489.15)
    , @(9,  0.00266,               11.75, 375.93, 1503.72)
    , @(10, 0.0026589999999999999, 2.9,   376.12, 1504.5)
    , @(11, 0.0052750000000000002, 47.39, 189.58, 758.32)
    , @(12, 0.0027160000000000001, 46.03, 368.23, 1472.9)
    , @(13, 0.0026770000000000001, 23.35, 373.55, 1494.22)
    , @(14, 0.0026679999999999998, 11.71, 374.87, 1499.47)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 7).Value  = $entry[1]   # G
    $ws.Cells.Item($row, 8).Value  = $entry[2]   # H
    $ws.Cells.Item($row, 9).Value  = $entry[3]   # I
    $ws.Cells.Item($row, 10).Value = $entry[4]   # J
}

# --- 4. Number formats -------------------------------------------------------
$ws.Range("G6:G14").NumberFormat = "0.000000"
$ws.Range("H6:J14").NumberFormat = "0.00"

# --- 5. Move the visible selection to K18 -----------------------------------
$ws.Range("K18").Select()

Write-Output "done"
